# Update countries & provincias Spain
# Reorders several country-name cells (column A) to reflect the new
# COVID-19 case-count sort order, refreshes the "datos actualizados" timestamp,
# and updates the numeric case/death columns (B:H) for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 00:20"
$ws.Range("B4").Value = 211463
$ws.Range("C4").Value = 22933
$ws.Range("E4").Value = 197937
$ws.Range("G4").Value = 668
$ws.Range("H4").Value = 4721
$ws.Range("C7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("B8").Value = 77921
$ws.Range("C8").Value = 6113
$ws.Range("E8").Value = 58296
$ws.Range("G8").Value = 150
$ws.Range("H8").Value = 925
$ws.Range("B18").Value = 9677
$ws.Range("C18").Value = 1065
$ws.Range("D18").Value = 1736
$ws.Range("E18").Value = 7827
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 114
$ws.Range("A22").Value = "Australia"
$ws.Range("B22").Value = 4980
$ws.Range("C22").Value = 217
$ws.Range("D22").Value = 345
$ws.Range("E22").Value = 4613
$ws.Range("F22").Value = 50
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 22
$ws.Range("A23").Value = "Suecia"
$ws.Range("B23").Value = 4947
$ws.Range("C23").Value = 512
$ws.Range("D23").Value = 103
$ws.Range("E23").Value = 4605
$ws.Range("F23").Value = 393
$ws.Range("G23").Value = 59
$ws.Range("H23").Value = 239
$ws.Range("B24").Value = 4877
$ws.Range("C24").Value = 236
$ws.Range("E24").Value = 4820
$ws.Range("D32").Value = 56
$ws.Range("E32").Value = 2455
$ws.Range("D47").Value = 236
$ws.Range("E47").Value = 982
$ws.Range("B119").Value = 101
$ws.Range("C119").Value = 7
$ws.Range("E119").Value = 90
$ws.Range("A125").Value = "Liechtenstein"
$ws.Range("B125").Value = 72
$ws.Range("D125").Value = 0
$ws.Range("E125").Value = 72
$ws.Range("F125").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("A126").Value = "Paraguay"
$ws.Range("B126").Value = 69
$ws.Range("C126").Value = 4
$ws.Range("D126").Value = 1
$ws.Range("E126").Value = 65
$ws.Range("F126").Value = 4
$ws.Range("H126").Value = 3
$ws.Range("A134").Value = "Jamaica"
$ws.Range("B134").Value = 44
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 2
$ws.Range("E134").Value = 39
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 3
$ws.Range("A135").Value = "Macao"
$ws.Range("B135").Value = 41
$ws.Range("D135").Value = 10
$ws.Range("E135").Value = 31
$ws.Range("H135").Value = 0
$ws.Range("A136").Value = "Puerto Rico"
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 1
$ws.Range("E136").Value = 36
$ws.Range("F136").Value = 0
$ws.Range("H136").Value = 2
$ws.Range("A137").Value = "Guatemala"
$ws.Range("B137").Value = 39
$ws.Range("C137").Value = 1
$ws.Range("D137").Value = 12
$ws.Range("E137").Value = 26
$ws.Range("F137").Value = 1
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 1
$ws.Range("A143").Value = "Guam"
$ws.Range("F143").Value = 0
$ws.Range("A144").Value = "El Salvador"
$ws.Range("E144").Value = 30
$ws.Range("F144").Value = 4
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 2
$ws.Range("A150").Value = "Congo"
$ws.Range("B150").Value = 22
$ws.Range("C150").Value = 3
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 20
$ws.Range("G150").Value = 2
$ws.Range("H150").Value = 2
$ws.Range("A151").Value = "Bahamas"
$ws.Range("B151").Value = 21
$ws.Range("C151").Value = 7
$ws.Range("E151").Value = 19
$ws.Range("G151").Value = 1
$ws.Range("A152").Value = "Tanzania"
$ws.Range("B152").Value = 20
$ws.Range("C152").Value = 1
$ws.Range("D152").Value = 1
$ws.Range("E152").Value = 18
$ws.Range("H152").Value = 1
$ws.Range("A153").Value = "Guyana"
$ws.Range("C153").Value = 7
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 16
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 3
$ws.Range("A154").Value = "Maldivas"
$ws.Range("B154").Value = 19
$ws.Range("C154").Value = 1
$ws.Range("D154").Value = 13
$ws.Range("E154").Value = 6
$ws.Range("H154").Value = 0
$ws.Range("A155").Value = "Gabon"
$ws.Range("B155").Value = 18
$ws.Range("C155").Value = 2
$ws.Range("H155").Value = 1
$ws.Range("A156").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B156").Value = 17
$ws.Range("E156").Value = 17
$ws.Range("A157").Value = "Nueva Caledonia"
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 0
$ws.Range("E157").Value = 16
$ws.Range("A158").Value = "Haiti"
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 1
$ws.Range("E158").Value = 15
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0
$ws.Range("A159").Value = "San Martin (Parte Holandesa)"
$ws.Range("B159").Value = 16
$ws.Range("C159").Value = 10
$ws.Range("D159").Value = 6
$ws.Range("E159").Value = 9
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 1
$ws.Range("A160").Value = "Eritrea"
$ws.Range("D160").Value = 0
$ws.Range("E160").Value = 15
$ws.Range("A162").Value = "Guinea Ecuatorial"
$ws.Range("D162").Value = 1
$ws.Range("E162").Value = 14
$ws.Range("H162").Value = 0
$ws.Range("A163").Value = "San Martin (Parte Francesa)"
$ws.Range("B163").Value = 15
$ws.Range("C163").Value = 0
$ws.Range("H163").Value = 1
$ws.Range("A164").Value = "Namibia"
$ws.Range("C164").Value = 3
$ws.Range("A165").Value = "Mongolia"
$ws.Range("B165").Value = 14
$ws.Range("C165").Value = 2
$ws.Range("D165").Value = 2
$ws.Range("A166").Value = "Santa Lucia"
$ws.Range("C166").Value = 0
$ws.Range("A167").Value = "Benin"
$ws.Range("B167").Value = 13
$ws.Range("C167").Value = 4
$ws.Range("D167").Value = 1
$ws.Range("A168").Value = "Dominica"
$ws.Range("E168").Value = 12
$ws.Range("H168").Value = 0
$ws.Range("A172").Value = "Laos"
$ws.Range("C172").Value = 1
$ws.Range("A173").Value = "Surinam"
$ws.Range("A174").Value = "Mozambique"
$ws.Range("C174").Value = 2
$ws.Range("A175").Value = "Groenlandia"
$ws.Range("D175").Value = 2
$ws.Range("H175").Value = 0
$ws.Range("A176").Value = "Siria"
$ws.Range("D176").Value = 0
$ws.Range("H176").Value = 2
$ws.Range("A177").Value = "Granada"
$ws.Range("A179").Value = "Suazilandia"
$ws.Range("A184").Value = "Republica del Chad"
$ws.Range("A185").Value = "Antigua y Barbuda"
$ws.Range("A188").Value = "Liberia"
$ws.Range("C188").Value = 3
$ws.Range("A189").Value = "Islas Turcas y Caicos"
$ws.Range("C189").Value = 1
$ws.Range("A190").Value = "San Bartolome"
$ws.Range("D190").Value = 1
$ws.Range("H190").Value = 0
$ws.Range("A191").Value = "Cabo Verde"
$ws.Range("D191").Value = 0
$ws.Range("H191").Value = 1
$ws.Range("A194").Value = "Somalia"
$ws.Range("A196").Value = "Nepal"
$ws.Range("A201").Value = "Belice"
$ws.Range("A202").Value = "Republica de Africa Central"
$ws.Range("A203").Value = "Sierra Leona"
$ws.Range("C203").Value = 1
$ws.Range("A205").Value = "Burundi"
$ws.Range("C205").Value = 0
$ws.Range("A206").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("C206").Value = 2
$ws.Range("A207").Value = "Papua Nueva Guinea"
$ws.Range("A208").Value = "Timor Oriental"
